$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1210.6666
$ws.Range("J17").Value = 1189.8948
$ws.Range("L17").Value = 3569.6844
$ws.Range("N17").Value = -3905.6844
$ws.Range("H19").Value = 2017.8148
$ws.Range("J19").Value = 1724.3572
$ws.Range("L19").Value = 1724.3572
$ws.Range("N19").Value = -2074.3572
$ws.Range("H33").Value = 797.04
$ws.Range("I33").Value = 884.5
$ws.Range("K33").Value = 884.5
$ws.Range("M33").Value = -655.5
$ws.Range("H37").Value = 3169.7144
$ws.Range("I37").Value = 1866.6666
$ws.Range("J37").Value = 4147
$ws.Range("K37").Value = 5599.9998
$ws.Range("L37").Value = 12441
$ws.Range("M37").Value = -5473.9998
$ws.Range("N37").Value = -12693
$ws.Range("H42").Value = 1101.8572
$ws.Range("J42").Value = 466.33334
$ws.Range("L42").Value = 1399.00002
$ws.Range("N42").Value = -1859.00002
$ws.Range("H46").Value = 874.75
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 899.5
$ws.Range("K46").Value = 2550
$ws.Range("L46").Value = 2698.5
$ws.Range("M46").Value = -2431
$ws.Range("N46").Value = -2936.5
$ws.Range("H60").Value = 874.75
$ws.Range("I60").Value = 850
$ws.Range("J60").Value = 899.5
$ws.Range("K60").Value = 2550
$ws.Range("L60").Value = 2698.5
$ws.Range("M60").Value = -2066
$ws.Range("N60").Value = -3666.5
$ws.Range("H62").Value = 3316.6667
$ws.Range("I62").Value = 2649.1428
$ws.Range("J62").Value = 5653
$ws.Range("K62").Value = 2649.1428
$ws.Range("L62").Value = 5653
$ws.Range("M62").Value = -2025.1428
$ws.Range("N62").Value = -6901
$ws.Range("H65").Value = 3316.6667
$ws.Range("I65").Value = 2649.1428
$ws.Range("J65").Value = 5653
$ws.Range("K65").Value = 13245.714
$ws.Range("L65").Value = 28265
$ws.Range("M65").Value = -10125.714
$ws.Range("N65").Value = -34505
$ws.Range("H70").Value = 1764.0454
$ws.Range("I70").Value = 1673.5
$ws.Range("K70").Value = 5020.5
$ws.Range("M70").Value = -4750.5
$ws.Range("H73").Value = 1764.0454
$ws.Range("I73").Value = 1673.5
$ws.Range("K73").Value = 5020.5
$ws.Range("M73").Value = -4084.5
$ws.Range("H80").Value = 1295.3846
$ws.Range("I80").Value = 1026.3
$ws.Range("J80").Value = 2192.3333
$ws.Range("K80").Value = 3078.9
$ws.Range("L80").Value = 6576.999899999999
$ws.Range("M80").Value = -2080.9
$ws.Range("N80").Value = -8572.999899999999
$ws.Range("H83").Value = 1295.3846
$ws.Range("I83").Value = 1026.3
$ws.Range("J83").Value = 2192.3333
$ws.Range("K83").Value = 9236.699999999999
$ws.Range("L83").Value = 19730.9997
$ws.Range("M83").Value = -4244.699999999999
$ws.Range("N83").Value = -29714.9997
$ws.Range("J86").Value = 2229.8
$ws.Range("L86").Value = 2229.8
$ws.Range("N86").Value = -4475.8
$ws.Range("J89").Value = 2229.8
$ws.Range("L89").Value = 11149
$ws.Range("N89").Value = -22381
$ws.Range("H98").Value = 1080.8
$ws.Range("I98").Value = 1080.8
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1080.8
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 417.2
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 1895.6842
$ws.Range("I100").Value = 1832.6666
$ws.Range("J100").Value = 1907.5
$ws.Range("K100").Value = 1832.6666
$ws.Range("L100").Value = 1907.5
$ws.Range("M100").Value = -1291.6666
$ws.Range("N100").Value = -2989.5
$ws.Range("H116").Value = 5939.174
$ws.Range("J116").Value = 4709.1816
$ws.Range("L116").Value = 4709.1816
$ws.Range("N116").Value = -11593.1816
$ws.Range("H122").Value = 1080.8
$ws.Range("I122").Value = 1080.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3242.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -792.3999999999996
$ws.Range("N122").ClearContents()
$ws.Range("H129").Value = 1129.6666
$ws.Range("I129").Value = 555.6
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 1666.8
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 3333.2
$ws.Range("N129").Value = -22000
$ws.Range("H131").Value = 3300.389
$ws.Range("I131").Value = 1386.2142
$ws.Range("K131").Value = 4158.642599999999
$ws.Range("M131").Value = 881.3574000000008
$ws.Range("H141").Value = 13549.417
$ws.Range("I141").Value = 17041.857
$ws.Range("J141").Value = 8660
$ws.Range("K141").Value = 51125.571
$ws.Range("L141").Value = 25980
$ws.Range("M141").Value = -45945.571
$ws.Range("N141").Value = -36340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2518.1145
$ws.Range("I32").Value = 2105.809
$ws.Range("J32").Value = 7760.2856
$ws.Range("K32").Value = 2105.809
$ws.Range("L32").Value = 7760.2856
$ws.Range("M32").Value = -1818.809
$ws.Range("N32").Value = -8334.285599999999
$ws.Range("H61").Value = 3306
$ws.Range("I61").Value = 3169.6667
$ws.Range("K61").Value = 3169.6667
$ws.Range("M61").Value = -2957.6667
$ws.Range("H74").Value = 4534.1494
$ws.Range("I74").Value = 4496.6875
$ws.Range("K74").Value = 4496.6875
$ws.Range("M74").Value = -3622.6875
$ws.Range("H77").Value = 4534.1494
$ws.Range("I77").Value = 4496.6875
$ws.Range("K77").Value = 22483.4375
$ws.Range("M77").Value = -18115.4375
$ws.Range("H81").Value = 9999
$ws.Range("I81").Value = 9999
$ws.Range("K81").Value = 9999
$ws.Range("M81").Value = -9001
$ws.Range("H84").Value = 9999
$ws.Range("I84").Value = 9999
$ws.Range("K84").Value = 29997
$ws.Range("M84").Value = -25005
$ws.Range("H102").Value = 4283.0713
$ws.Range("I102").Value = 4232
$ws.Range("K102").Value = 4232
$ws.Range("M102").Value = -2610
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H110").Value = 1981.3
$ws.Range("I110").Value = 2226.875
$ws.Range("J110").Value = 999
$ws.Range("K110").Value = 2226.875
$ws.Range("L110").Value = 999
$ws.Range("M110").Value = -181.875
$ws.Range("N110").Value = -5089
$ws.Range("H122").Value = 4678.1304
$ws.Range("I122").Value = 4401
$ws.Range("J122").Value = 4891.3076
$ws.Range("K122").Value = 13203
$ws.Range("L122").Value = 14673.9228
$ws.Range("M122").Value = -10753
$ws.Range("N122").Value = -19573.9228
$ws.Range("H132").Value = 3715.4082
$ws.Range("I132").Value = 3211.4243
$ws.Range("J132").Value = 4754.875
$ws.Range("K132").Value = 9634.2729
$ws.Range("L132").Value = 14264.625
$ws.Range("M132").Value = -7104.2729
$ws.Range("N132").Value = -19324.625
$ws.Range("H136").Value = 3306
$ws.Range("I136").Value = 3169.6667
$ws.Range("K136").Value = 9509.000100000001
$ws.Range("M136").Value = -6959.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2743.24
$ws.Range("I20").Value = 3221.3572
$ws.Range("J20").Value = 2134.7273
$ws.Range("K20").Value = 3221.3572
$ws.Range("L20").Value = 2134.7273
$ws.Range("M20").Value = -2974.3572
$ws.Range("N20").Value = -2628.7273
$ws.Range("H22").Value = 33601.8
$ws.Range("I22").Value = 252.07692
$ws.Range("J22").Value = 250375
$ws.Range("K22").Value = 252.07692
$ws.Range("L22").Value = 250375
$ws.Range("M22").Value = -79.07692
$ws.Range("N22").Value = -250721
$ws.Range("H30").Value = 1980
$ws.Range("J30").Value = 1980
$ws.Range("L30").Value = 1980
$ws.Range("N30").Value = -2230
$ws.Range("H35").Value = 100
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H50").Value = 69800
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H94").Value = 894.5
$ws.Range("I94").Value = 894.5
$ws.Range("K94").Value = 894.5
$ws.Range("M94").Value = -443.5
$ws.Range("H99").Value = 1967.1666
$ws.Range("I99").Value = 1879.5
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 1879.5
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = -381.5
$ws.Range("N99").Value = -5007
$ws.Range("H107").Value = 1858.92
$ws.Range("I107").Value = 2221
$ws.Range("J107").Value = 1688.5294
$ws.Range("K107").Value = 2221
$ws.Range("L107").Value = 1688.5294
$ws.Range("M107").Value = -301
$ws.Range("N107").Value = -5528.529399999999
$ws.Range("H130").Value = 89997.164
$ws.Range("J130").Value = 89997.164
$ws.Range("L130").Value = 89997.164
$ws.Range("N130").Value = -100037.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1589
$ws.Range("J16").Value = 1875.5555
$ws.Range("L16").Value = 1875.5555
$ws.Range("N16").Value = -2449.5555
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30460
$ws.Range("H31").Value = 2091.7222
$ws.Range("I31").Value = 1631.7778
$ws.Range("J31").Value = 2551.6667
$ws.Range("K31").Value = 1631.7778
$ws.Range("L31").Value = 2551.6667
$ws.Range("M31").Value = -1336.7778
$ws.Range("N31").Value = -3141.6667
$ws.Range("H34").Value = 2091.7222
$ws.Range("I34").Value = 1631.7778
$ws.Range("J34").Value = 2551.6667
$ws.Range("K34").Value = 1631.7778
$ws.Range("L34").Value = 2551.6667
$ws.Range("M34").Value = -1429.7778
$ws.Range("N34").Value = -2955.6667
$ws.Range("H39").Value = 1084443.8
$ws.Range("I39").Value = 1256876.1
$ws.Range("K39").Value = 1256876.1
$ws.Range("M39").Value = -1256485.1
$ws.Range("H49").Value = 1084443.8
$ws.Range("I49").Value = 1256876.1
$ws.Range("K49").Value = 1256876.1
$ws.Range("M49").Value = -1256694.1
$ws.Range("H58").Value = 11202.02
$ws.Range("I58").Value = 8138.3335
$ws.Range("K58").Value = 8138.3335
$ws.Range("M58").Value = -7935.3335
$ws.Range("H59").Value = 50052
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 100000
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 100000
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -102290
$ws.Range("H62").Value = 100003700
$ws.Range("J62").Value = 4609.4
$ws.Range("L62").Value = 4609.4
$ws.Range("N62").Value = -5857.4
$ws.Range("H65").Value = 100003700
$ws.Range("J65").Value = 4609.4
$ws.Range("L65").Value = 23047
$ws.Range("N65").Value = -29287
$ws.Range("H68").Value = 38770.57
$ws.Range("J68").Value = 38770.57
$ws.Range("L68").Value = 38770.57
$ws.Range("N68").Value = -40268.57
$ws.Range("H71").Value = 38770.57
$ws.Range("J71").Value = 38770.57
$ws.Range("L71").Value = 116311.71
$ws.Range("N71").Value = -123799.71
$ws.Range("H93").Value = 4324.8335
$ws.Range("I93").Value = 4516.8
$ws.Range("K93").Value = 4516.8
$ws.Range("M93").Value = -2644.8
$ws.Range("H105").Value = 3926.3076
$ws.Range("I105").Value = 1886.6666
$ws.Range("J105").Value = 5674.5713
$ws.Range("K105").Value = 1886.6666
$ws.Range("L105").Value = 5674.5713
$ws.Range("M105").Value = -139.6666
$ws.Range("N105").Value = -9168.5713
$ws.Range("H107").Value = 2313
$ws.Range("I107").Value = 1670
$ws.Range("J107").Value = 3277.5
$ws.Range("K107").Value = 1670
$ws.Range("L107").Value = 3277.5
$ws.Range("M107").Value = 250
$ws.Range("N107").Value = -7117.5
$ws.Range("H110").Value = 139999
$ws.Range("J110").Value = 139999
$ws.Range("L110").Value = 139999
$ws.Range("N110").Value = -148179
$ws.Range("H113").Value = 1589
$ws.Range("J113").Value = 1875.5555
$ws.Range("L113").Value = 1875.5555
$ws.Range("N113").Value = -6215.5555
$ws.Range("H121").Value = 64999
$ws.Range("J121").Value = 64999
$ws.Range("L121").Value = 64999
$ws.Range("N121").Value = -67619
$ws.Range("H122").Value = 4320.952
$ws.Range("I122").Value = 4076.625
$ws.Range("J122").Value = 4646.722
$ws.Range("K122").Value = 12229.875
$ws.Range("L122").Value = 13940.166
$ws.Range("M122").Value = -9779.875
$ws.Range("N122").Value = -18840.166
$ws.Range("H132").Value = 9436.4
$ws.Range("I132").Value = 3837.125
$ws.Range("J132").Value = 19390.666
$ws.Range("K132").Value = 11511.375
$ws.Range("L132").Value = 58171.99800000001
$ws.Range("M132").Value = -8981.375
$ws.Range("N132").Value = -63231.99800000001
$ws.Range("H136").Value = 11202.02
$ws.Range("I136").Value = 8138.3335
$ws.Range("K136").Value = 24415.0005
$ws.Range("M136").Value = -21865.0005
$ws.Range("H138").Value = 103218.586
$ws.Range("J138").Value = 103218.586
$ws.Range("L138").Value = 103218.586
$ws.Range("N138").Value = -113498.586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 210.23077
$ws.Range("I6").Value = 159.11111
$ws.Range("J6").Value = 325.25
$ws.Range("K6").Value = 477.33333
$ws.Range("L6").Value = 975.75
$ws.Range("M6").Value = -364.33333
$ws.Range("N6").Value = -1201.75
$ws.Range("H13").Value = 4057.5454
$ws.Range("I13").Value = 399.75
$ws.Range("K13").Value = 1199.25
$ws.Range("M13").Value = -1031.25
$ws.Range("H34").Value = 1270
$ws.Range("I34").Value = 1580
$ws.Range("J34").Value = 495
$ws.Range("K34").Value = 4740
$ws.Range("L34").Value = 1485
$ws.Range("M34").Value = -4656
$ws.Range("N34").Value = -1653
$ws.Range("H39").Value = 6985.1816
$ws.Range("I39").Value = 9000
$ws.Range("J39").Value = 6537.4443
$ws.Range("K39").Value = 27000
$ws.Range("L39").Value = 19612.3329
$ws.Range("M39").Value = -26706
$ws.Range("N39").Value = -20200.3329
$ws.Range("H55").Value = 1104.5217
$ws.Range("J55").Value = 1050.1818
$ws.Range("L55").Value = 3150.5454
$ws.Range("N55").Value = -3504.5454
$ws.Range("H57").Value = 5859.7
$ws.Range("J57").Value = 5324.625
$ws.Range("L57").Value = 15973.875
$ws.Range("N57").Value = -17091.875
$ws.Range("H113").Value = 2029.1052
$ws.Range("I113").Value = 2097.25
$ws.Range("J113").Value = 2010.9333
$ws.Range("K113").Value = 6291.75
$ws.Range("L113").Value = 6032.7999
$ws.Range("M113").Value = -4121.75
$ws.Range("N113").Value = -10372.7999
$ws.Range("H119").Value = 5452.8887
$ws.Range("I119").Value = 4173.143
$ws.Range("K119").Value = 12519.429
$ws.Range("M119").Value = -7681.429
$ws.Range("H129").Value = 1677.75
$ws.Range("I129").Value = 646.7692
$ws.Range("J129").Value = 3592.4285
$ws.Range("K129").Value = 1940.3076
$ws.Range("L129").Value = 10777.2855
$ws.Range("M129").Value = 3059.6924
$ws.Range("N129").Value = -20777.2855
$ws.Range("H132").Value = 3940.7144
$ws.Range("I132").Value = 3223.25
$ws.Range("J132").Value = 4897.3335
$ws.Range("K132").Value = 29009.25
$ws.Range("L132").Value = 44076.0015
$ws.Range("M132").Value = -26479.25
$ws.Range("N132").Value = -49136.0015
$ws.Range("H134").Value = 2236.4707
$ws.Range("I134").Value = 1232.3077
$ws.Range("K134").Value = 3696.9231
$ws.Range("M134").Value = 1373.0769
$ws.Range("H139").Value = 3354.16
$ws.Range("I139").Value = 866.8421
$ws.Range("K139").Value = 2600.5263
$ws.Range("M139").Value = 2539.4737
$ws.Range("H140").Value = 1869.625
$ws.Range("I140").Value = 1869.625
$ws.Range("K140").Value = 5608.875
$ws.Range("M140").Value = -428.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 11289.667
$ws.Range("I99").Value = 6060.125
$ws.Range("J99").Value = 21748.75
$ws.Range("K99").Value = 6060.125
$ws.Range("L99").Value = 21748.75
$ws.Range("M99").Value = -3814.125
$ws.Range("N99").Value = -26240.75
$ws.Range("H110").Value = 145567.33
$ws.Range("J110").Value = 145567.33
$ws.Range("L110").Value = 145567.33
$ws.Range("N110").Value = -153747.33
$ws.Range("H113").Value = 2750.1667
$ws.Range("I113").Value = 1493.5
$ws.Range("J113").Value = 5263.5
$ws.Range("K113").Value = 1493.5
$ws.Range("L113").Value = 5263.5
$ws.Range("M113").Value = 676.5
$ws.Range("N113").Value = -9603.5
$ws.Range("H122").Value = 3635.5789
$ws.Range("I122").Value = 2828.1177
$ws.Range("K122").Value = 8484.3531
$ws.Range("M122").Value = -6034.3531
$ws.Range("H132").Value = 6919.2104
$ws.Range("I132").Value = 6812.273
$ws.Range("J132").Value = 7066.25
$ws.Range("K132").Value = 20436.819
$ws.Range("L132").Value = 21198.75
$ws.Range("M132").Value = -17906.819
$ws.Range("N132").Value = -26258.75
$ws.Range("H135").Value = 156726.33
$ws.Range("J135").Value = 157603.5
$ws.Range("L135").Value = 157603.5
$ws.Range("N135").Value = -167743.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2394
$ws.Range("I7").Value = 1999.5
$ws.Range("J7").Value = 2788.5
$ws.Range("K7").Value = 1999.5
$ws.Range("L7").Value = 2788.5
$ws.Range("M7").Value = -1887.5
$ws.Range("N7").Value = -3012.5
$ws.Range("H16").Value = 4035184.5
$ws.Range("I16").Value = 4631734
$ws.Range("K16").Value = 4631734
$ws.Range("M16").Value = -4631564
$ws.Range("H22").Value = 865.3333
$ws.Range("I22").Value = 863.125
$ws.Range("J22").Value = 869.75
$ws.Range("K22").Value = 863.125
$ws.Range("L22").Value = 869.75
$ws.Range("M22").Value = -568.125
$ws.Range("N22").Value = -1459.75
$ws.Range("H27").Value = 865.3333
$ws.Range("I27").Value = 863.125
$ws.Range("J27").Value = 869.75
$ws.Range("K27").Value = 863.125
$ws.Range("L27").Value = 869.75
$ws.Range("M27").Value = -756.125
$ws.Range("N27").Value = -1083.75
$ws.Range("H40").Value = 4973.9062
$ws.Range("I40").Value = 4952.2173
$ws.Range("J40").Value = 5029.3335
$ws.Range("K40").Value = 4952.2173
$ws.Range("L40").Value = 5029.3335
$ws.Range("M40").Value = -4816.2173
$ws.Range("N40").Value = -5301.3335
$ws.Range("H46").Value = 1939.3438
$ws.Range("I46").Value = 1462.2106
$ws.Range("J46").Value = 2636.6924
$ws.Range("K46").Value = 1462.2106
$ws.Range("L46").Value = 2636.6924
$ws.Range("M46").Value = -1274.2106
$ws.Range("N46").Value = -3012.6924
$ws.Range("H97").Value = 51748.75
$ws.Range("J97").Value = 51748.75
$ws.Range("L97").Value = 51748.75
$ws.Range("N97").Value = -53730.75
$ws.Range("H100").Value = 111113120
$ws.Range("I100").Value = 142859000
$ws.Range("J100").Value = 2499
$ws.Range("K100").Value = 142859000
$ws.Range("L100").Value = 2499
$ws.Range("M100").Value = -142858459
$ws.Range("N100").Value = -3581
$ws.Range("H126").Value = 2394
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 2788.5
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 8365.5
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -13305.5
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H130").Value = 44958.332
$ws.Range("J130").Value = 44958.332
$ws.Range("L130").Value = 44958.332
$ws.Range("N130").Value = -54998.332
$ws.Range("H132").Value = 4861.404
$ws.Range("I132").Value = 4614.7075
$ws.Range("K132").Value = 13844.1225
$ws.Range("M132").Value = -11314.1225

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 83338370
$ws.Range("I62").Value = 125004984
$ws.Range("J62").Value = 5125
$ws.Range("K62").Value = 125004984
$ws.Range("L62").Value = 5125
$ws.Range("M62").Value = -125004360
$ws.Range("N62").Value = -6373
$ws.Range("H65").Value = 83338370
$ws.Range("I65").Value = 125004984
$ws.Range("J65").Value = 5125
$ws.Range("K65").Value = 625024920
$ws.Range("L65").Value = 25625
$ws.Range("M65").Value = -625021800
$ws.Range("N65").Value = -31865
$ws.Range("H81").Value = 3034221.8
$ws.Range("I81").Value = 3640126.2
$ws.Range("J81").Value = 4699
$ws.Range("K81").Value = 7280252.4
$ws.Range("L81").Value = 9398
$ws.Range("M81").Value = -7279191.4
$ws.Range("N81").Value = -11520
$ws.Range("H84").Value = 3034221.8
$ws.Range("I84").Value = 3640126.2
$ws.Range("J84").Value = 4699
$ws.Range("K84").Value = 36401262
$ws.Range("L84").Value = 46990
$ws.Range("M84").Value = -36395958
$ws.Range("N84").Value = -57598
$ws.Range("H107").Value = 7428.7144
$ws.Range("I107").Value = 6999.6665
$ws.Range("K107").Value = 20998.9995
$ws.Range("M107").Value = -19078.9995
$ws.Range("H122").Value = 5432.276
$ws.Range("I122").Value = 2668.5881
$ws.Range("J122").Value = 9347.5
$ws.Range("K122").Value = 8005.7643
$ws.Range("L122").Value = 28042.5
$ws.Range("M122").Value = -5555.7643
$ws.Range("N122").Value = -32942.5
$ws.Range("H128").Value = 99997.5
$ws.Range("J128").Value = 99997.5
$ws.Range("L128").Value = 99997.5
$ws.Range("N128").Value = -109957.5
$ws.Range("H132").Value = 2699.8604
$ws.Range("I132").Value = 2168.9048
$ws.Range("J132").Value = 25000
$ws.Range("K132").Value = 6506.714399999999
$ws.Range("L132").Value = 75000
$ws.Range("M132").Value = -3976.714399999999
$ws.Range("N132").Value = -80060
$ws.Range("H137").Value = 114960.445
$ws.Range("J137").Value = 114960.445
$ws.Range("L137").Value = 114960.445
$ws.Range("N137").Value = -125160.445
